$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically swaps the two blocks of 5 data rows:
#   rows 2-6  (old)  <->  rows 7-11 (old)
# i.e. new row 2 = old row 7, new row 3 = old row 8, ..., new row 6 = old row 11
#      new row 7 = old row 2, new row 8 = old row 3, ..., new row 11 = old row 6
# Columns span A (1) through AY (51).

$firstCol = 1
$lastCol = 51
$blockSize = 5
$topStart = 2       # rows 2..6
$bottomStart = 7     # rows 7..11

# Columns that are stored as text in this sheet (vs. numeric/boolean). Some of
# their values look like numbers or ISO dates (e.g. "2", "2023-07-10"), and a
# plain .Value assignment would let Excel auto-coerce them into numbers/dates.
# Forcing the destination cell to Text format before the write keeps them as
# text, matching the original column typing.
$textCols = @(3,4,6,7,8,9,11,12,13,14,16,20,21,22,23,25,26,27,28,45,46,49,50,51)

for ($c = $firstCol; $c -le $lastCol; $c++) {
    $isTextCol = $textCols -contains $c
    for ($i = 0; $i -lt $blockSize; $i++) {
        $topRow = $topStart + $i
        $bottomRow = $bottomStart + $i

        $topCell = $ws.Cells.Item($topRow, $c)
        $bottomCell = $ws.Cells.Item($bottomRow, $c)

        $topVal = $topCell.Value()
        $bottomVal = $bottomCell.Value()

        if ($isTextCol) {
            $topCell.NumberFormat = "@"
            $bottomCell.NumberFormat = "@"
        }

        $topCell.Value = $bottomVal
        $bottomCell.Value = $topVal
    }
}
